$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.5547, 1.091, 0, 15),
    @(3, 0.7871, 1.041, 0, 20),
    @(4, 0.6943, 1.125, 0, 20),
    @(5, 0.6021, 1.216, 0, 20),
    @(6, 0.4131, 1.253, 0, 15),
    @(7, 0.351, 1.352, 0, 15),
    @(8, 1.092, 1.095, 0, 25),
    @(9, 0.1939, 1.63, 0, 10),
    @(10, 0.6674, 1.492, 0, 20),
    @(11, 0.4228, 1.817, 0, 15),
    @(12, 0.7419, 1.594, 0, 20),
    @(13, 0.7859, 1.654, 0, 20),
    @(14, 0.8205, 1.797, 0, 20),
    @(15, 1.012, 1.534, 0, 25),
    @(16, 0.1007, 2.327, 0, 10),
    @(17, 0.00008296000000000001, 3.598, 0.3733, 5),
    @(18, 0.02765, 2.838, 0, 10),
    @(19, 0.1738, 2.872, 0, 15),
    @(20, 0.03501, 3.167, 0, 10),
    @(21, 0.1341, 3.207, 0, 15),
    @(22, 0.05612, 3.342, 0, 15),
    @(23, 0.262, 3.172, 0, 20),
    @(24, 0.1853, 3.307, 0, 20),
    @(25, 0.5451, 2.916, 0, 25),
    @(26, 0.0008876999999999999, 5.263, 0.1375, 5),
    @(27, 0.02503, 4.589, 0, 10),
    @(28, 0.01462, 4.589, 0, 10),
    @(29, 0.0077, 4.657, 0, 10),
    @(30, 0.0001061, 5.513, 0.07238, 5),
    @(31, 0.06807000000000001, 4.737, 0, 15),
    @(32, 0.143, 4.348, 0, 20),
    @(33, 0.1074, 4.427, 0, 20),
    @(34, 0.4929, 3.961, 0, 25),
    @(35, 0.3971, 4.047, 0, 25),
    @(36, 0.005569, 4.771, 0, 10),
    @(37, 0.07881000000000001, 4.998, 0, 15),
    @(38, 0.001069, 6.106, 0.1564, 5),
    @(39, 0.2017, 4.984, 0, 20),
    @(40, 0.1572, 4.968, 0, 20),
    @(41, 0.4605, 4.586, 0, 25),
    @(42, 0.006423, 5.057, 0, 10),
    @(43, 0.1434, 5.187, 0, 15),
    @(44, 0.6493, 4.823, 0, 25),
    @(45, 0.6409, 4.865, 0, 25),
    @(46, 0.1509, 5.343, 0, 10),
    @(47, 0.5592, 4.956, 0, 25),
    @(48, 0.406, 5.559, 0, 15),
    @(49, 0.6852, 5.047, 0, 25),
    @(50, 1.08, 5.027, 0, 25),
    @(51, 0.003351, 6.489, 0.6151, 5),
    @(52, 0.04785, 5.734, 0, 10),
    @(53, 0.05953, 5.575, 0, 10),
    @(54, 0.6843, 5.541, 0, 15),
    @(55, 1.42, 5.093, 0, 25),
    @(56, 1.354, 5.316, 0, 20),
    @(57, 1.348, 5.283, 0, 20),
    @(58, 1.832, 5.019, 0, 25),
    @(59, 0.3051, 5.137, 0, 15),
    @(60, 1.797, 4.918, 0, 25),
    @(61, 0.194, 4.895, 0, 15),
    @(62, 0.147, 4.801, 0, 15),
    @(63, 0.216, 4.72, 0, 15),
    @(64, 0.1372, 4.639, 0, 15),
    @(65, 0.004142, 5.546, 0.2186, 5),
    @(66, 0.01208, 4.676, 0, 10),
    @(67, 0.01316, 4.738, 0, 10),
    @(68, 0.1508, 4.885, 0, 15),
    @(69, 0.1178, 5.017, 0, 15),
    @(70, 0.02316, 5.41, 0, 10),
    @(71, 0.3763, 5.059, 0, 20),
    @(72, 0.02234, 5.85, 0, 10)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}